$d = $word.ActiveDocument

# Locate the paragraph that marks the start of the EvolutionIQ job's date/location
# line ("New York City, NY - Feb 2024 - Current"). We scan Paragraphs rather than
# relying on a fixed index so the script stays correct if content shifts.
$anchorIndex = -1
$i = 0
foreach ($para in $d.Paragraphs) {
    $i = $i + 1
    if ($para.Range.Text -like "New York City, NY - Feb 2024 - Current*") {
        $anchorIndex = $i
    }
}

if ($anchorIndex -eq -1) {
    throw "Could not find anchor paragraph 'New York City, NY - Feb 2024 - Current'"
}

$anchorPara = $d.Paragraphs.Item($anchorIndex)
$nextPara = $anchorPara.Next()

# Insert a brand-new paragraph immediately before the next paragraph (the first
# EvolutionIQ bullet). It inherits that bullet's paragraph formatting (Normal
# style + numPr ilvl=0/numId=3 bullet list), matching the target formatting.
$insertionPoint = $nextPara.Range
$insertionPoint.Collapse(1)
$insertionPoint.InsertParagraphBefore()

$newPara = $anchorPara.Next()
$newPara.Range.InsertBefore("Created infrastructure to provide ephemeral environments, allowing testing of code branches before sending those code changes to the company at large.")
